$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.9568136666666667
$ws.Cells.Item(2, 8).Value = 2.870441
$ws.Cells.Item(2, 9).Value = 0.09967139189263423
$ws.Cells.Item(2, 10).Value = 0.09967139189263421
$ws.Cells.Item(2, 13).Value = 3.759736666666667
$ws.Cells.Item(2, 14).Value = 11.27921
$ws.Cells.Item(2, 15).Value = 0.0683751702595819
$ws.Cells.Item(2, 16).Value = 0.06837517025958188
$ws.Cells.Item(2, 17).Value = 3.597367425734444
$ws.Cells.Item(2, 18).Value = 32.37630683160999
$ws.Cells.Item(2, 19).Value = 0.006815048390668376
$ws.Cells.Item(2, 20).Value = 0.006815048390668374

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.9568136666666667
$ws.Cells.Item(3, 8).Value = 2.870441
$ws.Cells.Item(3, 9).Value = 0.09967139189263423
$ws.Cells.Item(3, 10).Value = 0.09967139189263421
$ws.Cells.Item(3, 15).Value = 0.6514180024294648
$ws.Cells.Item(3, 16).Value = 0.6514180024294647
$ws.Cells.Item(3, 17).Value = 34.27252749178145
$ws.Cells.Item(3, 18).Value = 308.452747426033
$ws.Cells.Item(3, 19).Value = 0.06492773900606415
$ws.Cells.Item(3, 20).Value = 0.06492773900606412

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.9568136666666667
$ws.Cells.Item(4, 8).Value = 2.870441
$ws.Cells.Item(4, 9).Value = 0.09967139189263423
$ws.Cells.Item(4, 10).Value = 0.09967139189263421
$ws.Cells.Item(4, 15).Value = 0.2802068273109533
$ws.Cells.Item(4, 16).Value = 0.2802068273109533
$ws.Cells.Item(4, 17).Value = 14.74229474252111
$ws.Cells.Item(4, 18).Value = 132.68065268269
$ws.Cells.Item(4, 19).Value = 0.02792860449590171
$ws.Cells.Item(4, 20).Value = 0.0279286044959017

# Row 5
$ws.Cells.Item(5, 9).Value = 0.789130862182032
$ws.Cells.Item(5, 10).Value = 0.789130862182032
$ws.Cells.Item(5, 13).Value = 3.759736666666667
$ws.Cells.Item(5, 14).Value = 11.27921
$ws.Cells.Item(5, 15).Value = 0.0683751702595819
$ws.Cells.Item(5, 16).Value = 0.06837517025958188
$ws.Cells.Item(5, 17).Value = 28.48152919659556
$ws.Cells.Item(5, 18).Value = 256.33376276936
$ws.Cells.Item(5, 19).Value = 0.05395695705878709
$ws.Cells.Item(5, 20).Value = 0.05395695705878709

# Row 6
$ws.Cells.Item(6, 9).Value = 0.789130862182032
$ws.Cells.Item(6, 10).Value = 0.789130862182032
$ws.Cells.Item(6, 15).Value = 0.6514180024294648
$ws.Cells.Item(6, 16).Value = 0.6514180024294647
$ws.Cells.Item(6, 19).Value = 0.5140540498980606
$ws.Cells.Item(6, 20).Value = 0.5140540498980605

# Row 7
$ws.Cells.Item(7, 9).Value = 0.789130862182032
$ws.Cells.Item(7, 10).Value = 0.789130862182032
$ws.Cells.Item(7, 15).Value = 0.2802068273109533
$ws.Cells.Item(7, 16).Value = 0.2802068273109533
$ws.Cells.Item(7, 19).Value = 0.2211198552251843
$ws.Cells.Item(7, 20).Value = 0.2211198552251843

# Row 8
$ws.Cells.Item(8, 9).Value = 0.1111977459253338
$ws.Cells.Item(8, 10).Value = 0.1111977459253338
$ws.Cells.Item(8, 13).Value = 3.759736666666667
$ws.Cells.Item(8, 14).Value = 11.27921
$ws.Cells.Item(8, 15).Value = 0.0683751702595819
$ws.Cells.Item(8, 16).Value = 0.06837517025958188
$ws.Cells.Item(8, 17).Value = 4.01337978141
$ws.Cells.Item(8, 18).Value = 36.12041803269
$ws.Cells.Item(8, 19).Value = 0.007603164810126427
$ws.Cells.Item(8, 20).Value = 0.007603164810126425

# Row 9
$ws.Cells.Item(9, 9).Value = 0.1111977459253338
$ws.Cells.Item(9, 10).Value = 0.1111977459253338
$ws.Cells.Item(9, 15).Value = 0.6514180024294648
$ws.Cells.Item(9, 16).Value = 0.6514180024294647
$ws.Cells.Item(9, 19).Value = 0.07243621352534009
$ws.Cells.Item(9, 20).Value = 0.07243621352534008

# Row 10
$ws.Cells.Item(10, 9).Value = 0.1111977459253338
$ws.Cells.Item(10, 10).Value = 0.1111977459253338
$ws.Cells.Item(10, 15).Value = 0.2802068273109533
$ws.Cells.Item(10, 16).Value = 0.2802068273109533
$ws.Cells.Item(10, 19).Value = 0.03115836758986727
$ws.Cells.Item(10, 20).Value = 0.03115836758986726


$wb.Save()
